# Auto-generated edit script applying cell-level data updates
# (market price refresh) per the authoritative diff.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2033.4546
$ws.Range("I6").Value = 171
$ws.Range("K6").Value = 513
$ws.Range("M6").Value = -401
$ws.Range("H53").Value = 43823.566
$ws.Range("J53").Value = 91243.27
$ws.Range("L53").Value = 91243.27
$ws.Range("N53").Value = -92517.27
$ws.Range("H76").Value = 10889167
$ws.Range("I76").Value = 355999.8
$ws.Range("J76").Value = 28444444
$ws.Range("K76").Value = 355999.8
$ws.Range("L76").Value = 28444444
$ws.Range("M76").Value = -355684.8
$ws.Range("N76").Value = -28445074
$ws.Range("H79").Value = 10889167
$ws.Range("I79").Value = 355999.8
$ws.Range("J79").Value = 28444444
$ws.Range("K79").Value = 355999.8
$ws.Range("L79").Value = 28444444
$ws.Range("M79").Value = -354907.8
$ws.Range("N79").Value = -28446628
$ws.Range("H86").Value = 7734148
$ws.Range("J86").Value = 25129140
$ws.Range("L86").Value = 25129140
$ws.Range("N86").Value = -25131386
$ws.Range("H88").Value = 3140.7856
$ws.Range("J88").Value = 3287.7
$ws.Range("L88").Value = 3287.7
$ws.Range("N88").Value = -4099.7
$ws.Range("H89").Value = 7734148
$ws.Range("J89").Value = 25129140
$ws.Range("L89").Value = 125645700
$ws.Range("N89").Value = -125656932
$ws.Range("H91").Value = 3140.7856
$ws.Range("J91").Value = 3287.7
$ws.Range("L91").Value = 3287.7
$ws.Range("N91").Value = -6095.7
$ws.Range("H98").Value = 577.0476
$ws.Range("I98").Value = 503.75674
$ws.Range("K98").Value = 503.75674
$ws.Range("M98").Value = 994.24326
$ws.Range("H99").Value = 486
$ws.Range("I99").Value = 486
$ws.Range("K99").Value = 1458
$ws.Range("M99").Value = 40
$ws.Range("H104").Value = 3188.8333
$ws.Range("I104").Value = 3504.75
$ws.Range("J104").Value = 2557
$ws.Range("K104").Value = 10514.25
$ws.Range("L104").Value = 7671
$ws.Range("M104").Value = -8767.25
$ws.Range("N104").Value = -11165
$ws.Range("H112").Value = 4253.857
$ws.Range("I112").Value = 4755.8
$ws.Range("J112").Value = 2999
$ws.Range("K112").Value = 14267.4
$ws.Range("L112").Value = 8997
$ws.Range("M112").Value = -13159.4
$ws.Range("N112").Value = -11213
$ws.Range("H122").Value = 577.0476
$ws.Range("I122").Value = 503.75674
$ws.Range("K122").Value = 1511.27022
$ws.Range("M122").Value = 938.7297800000001
$ws.Range("H131").Value = 9630.556
$ws.Range("I131").Value = 2968.75
$ws.Range("K131").Value = 8906.25
$ws.Range("M131").Value = -3866.25
$ws.Range("H132").Value = 2492.9453
$ws.Range("I132").Value = 2150.2285
$ws.Range("J132").Value = 10489.667
$ws.Range("K132").Value = 6450.685500000001
$ws.Range("L132").Value = 31469.001
$ws.Range("M132").Value = -3920.685500000001
$ws.Range("N132").Value = -36529.001
$ws.Range("H137").Value = 2200
$ws.Range("I137").Value = 2200
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6600
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -4050
$ws.Range("N137").ClearContents()
$ws.Range("H141").Value = 850.8570999999999
$ws.Range("I141").Value = 882.73334
$ws.Range("K141").Value = 2648.20002
$ws.Range("M141").Value = 2531.79998

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7546.7407
$ws.Range("I32").Value = 6150.52
$ws.Range("K32").Value = 6150.52
$ws.Range("M32").Value = -5863.52
$ws.Range("H97").Value = 2008.5385
$ws.Range("I97").Value = 1752.125
$ws.Range("J97").Value = 2418.8
$ws.Range("K97").Value = 1752.125
$ws.Range("L97").Value = 2418.8
$ws.Range("M97").Value = -1256.125
$ws.Range("N97").Value = -3410.8
$ws.Range("H102").Value = 40362.5
$ws.Range("I102").Value = 3238.2354
$ws.Range("K102").Value = 3238.2354
$ws.Range("M102").Value = -1616.2354

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 26102
$ws.Range("J20").Value = 26431
$ws.Range("L20").Value = 26431
$ws.Range("N20").Value = -26925
$ws.Range("H99").Value = 3634.5
$ws.Range("I99").Value = 4235.222
$ws.Range("K99").Value = 4235.222
$ws.Range("M99").Value = -2737.222
$ws.Range("H105").Value = 1746.909
$ws.Range("I105").Value = 1569.0555
$ws.Range("J105").Value = 2547.25
$ws.Range("K105").Value = 1569.0555
$ws.Range("L105").Value = 2547.25
$ws.Range("M105").Value = 177.9445000000001
$ws.Range("N105").Value = -6041.25
$ws.Range("H107").Value = 2076.7083
$ws.Range("I107").Value = 1665.3684
$ws.Range("K107").Value = 1665.3684
$ws.Range("M107").Value = 254.6315999999999
$ws.Range("H134").Value = 930.7727
$ws.Range("I134").Value = 874.15
$ws.Range("K134").Value = 2622.45
$ws.Range("M134").Value = -87.44999999999982
$ws.Range("H135").Value = 41375
$ws.Range("J135").Value = 41375
$ws.Range("L135").Value = 41375
$ws.Range("N135").Value = -51515

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5977.893
$ws.Range("I31").Value = 7554
$ws.Range("K31").Value = 7554
$ws.Range("M31").Value = -7259
$ws.Range("H34").Value = 5977.893
$ws.Range("I34").Value = 7554
$ws.Range("K34").Value = 7554
$ws.Range("M34").Value = -7352
$ws.Range("H58").Value = 1583.963
$ws.Range("I58").Value = 948.6957
$ws.Range("J58").Value = 5236.75
$ws.Range("K58").Value = 948.6957
$ws.Range("L58").Value = 5236.75
$ws.Range("M58").Value = -745.6957
$ws.Range("N58").Value = -5642.75
$ws.Range("H134").Value = 1566.2812
$ws.Range("I134").Value = 1462.4348
$ws.Range("J134").Value = 1831.6666
$ws.Range("K134").Value = 4387.3044
$ws.Range("L134").Value = 5494.9998
$ws.Range("M134").Value = -1852.3044
$ws.Range("N134").Value = -10564.9998
$ws.Range("H136").Value = 1583.963
$ws.Range("I136").Value = 948.6957
$ws.Range("J136").Value = 5236.75
$ws.Range("K136").Value = 2846.0871
$ws.Range("L136").Value = 15710.25
$ws.Range("M136").Value = -296.0870999999997
$ws.Range("N136").Value = -20810.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 585
$ws.Range("I14").Value = 585
$ws.Range("K14").Value = 1755
$ws.Range("M14").Value = -1582
$ws.Range("H107").Value = 2250
$ws.Range("J107").Value = 1250
$ws.Range("L107").Value = 3750
$ws.Range("N107").Value = -7590

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H70").Value = 13354.909
$ws.Range("I70").Value = 13666.667
$ws.Range("J70").Value = 12980.8
$ws.Range("K70").Value = 13666.667
$ws.Range("L70").Value = 12980.8
$ws.Range("M70").Value = -13396.667
$ws.Range("N70").Value = -13520.8
$ws.Range("H73").Value = 13354.909
$ws.Range("I73").Value = 13666.667
$ws.Range("J73").Value = 12980.8
$ws.Range("K73").Value = 13666.667
$ws.Range("L73").Value = 12980.8
$ws.Range("M73").Value = -12730.667
$ws.Range("N73").Value = -14852.8
$ws.Range("H119").Value = 55330
$ws.Range("J119").Value = 55330
$ws.Range("L119").Value = 55330
$ws.Range("N119").Value = -65006
$ws.Range("H132").Value = 3516
$ws.Range("I132").Value = 2913.375
$ws.Range("J132").Value = 10747.5
$ws.Range("K132").Value = 8740.125
$ws.Range("L132").Value = 32242.5
$ws.Range("M132").Value = -6210.125
$ws.Range("N132").Value = -37302.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1796.0938
$ws.Range("I46").Value = 1585.238
$ws.Range("J46").Value = 2198.6365
$ws.Range("K46").Value = 1585.238
$ws.Range("L46").Value = 2198.6365
$ws.Range("M46").Value = -1397.238
$ws.Range("N46").Value = -2574.6365
$ws.Range("H55").Value = 444.35715
$ws.Range("I55").Value = 438.05
$ws.Range("J55").Value = 460.125
$ws.Range("K55").Value = 438.05
$ws.Range("L55").Value = 460.125
$ws.Range("M55").Value = -265.05
$ws.Range("N55").Value = -806.125
$ws.Range("H136").Value = 2070.7693
$ws.Range("I136").Value = 1738.6061
$ws.Range("J136").Value = 3897.6667
$ws.Range("K136").Value = 5215.8183
$ws.Range("L136").Value = 11693.0001
$ws.Range("M136").Value = -2665.8183
$ws.Range("N136").Value = -16793.0001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 80309120
$ws.Range("I62").Value = 4872891.5
$ws.Range("J62").Value = 250040620
$ws.Range("K62").Value = 4872891.5
$ws.Range("L62").Value = 250040620
$ws.Range("M62").Value = -4872267.5
$ws.Range("N62").Value = -250041868
$ws.Range("H65").Value = 80309120
$ws.Range("I65").Value = 4872891.5
$ws.Range("J65").Value = 250040620
$ws.Range("K65").Value = 24364457.5
$ws.Range("L65").Value = 1250203100
$ws.Range("M65").Value = -24361337.5
$ws.Range("N65").Value = -1250209340
$ws.Range("H107").Value = 916.2308
$ws.Range("J107").Value = 1399.3334
$ws.Range("L107").Value = 4198.0002
$ws.Range("N107").Value = -8038.0002
$ws.Range("H136").Value = 3304.7222
$ws.Range("I136").Value = 2216.25
$ws.Range("J136").Value = 5481.6665
$ws.Range("K136").Value = 6648.75
$ws.Range("L136").Value = 16444.9995
$ws.Range("M136").Value = -4098.75
$ws.Range("N136").Value = -21544.9995

